$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6 from 2023-10-22 (45221)
# to 2023-10-25 (45224), preserving the existing date number format on these cells.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45224
}
